# Sprint 3 / Sprint 4 schedule date corrections in the "Sprint" summary
# table (Sprint, Total Story Points, Duration, Sprint Start Date,
# Sprint End Date (Planned), Story Points Completed, Sprint Release
# Date (Actual)).
#
#   Sprint-3 row: Start Date   07 Nov 2022 -> 06 Nov 2022
#                 End Date     12 Nov 2022 -> 11 Nov 2022
#                 Release Date 12 Nov 2022 -> 09 Nov 2022
#   Sprint-4 row: Start Date   14 Nov 2022 -> 10 Nov 2022
#                 End Date     19 Nov 2022 -> 15 Nov 2022
#                 Release Date 19 Nov 2022 -> 15 Nov 2022

$d = $word.ActiveDocument

# The sprint-summary table is the 3rd table in the document.
$sprintTable = $d.Tables.Item(3)

function Set-CellText {
    param(
        $Table,
        [int]$Row,
        [int]$Col,
        [string]$NewText
    )

    $cell = $Table.Cell($Row, $Col)
    $rng = $cell.Range
    # Exclude the trailing cell-mark character from the range so the
    # whole visible run is replaced in one shot (keeps the run's
    # original formatting / rsid intact instead of fragmenting it).
    $rng.End = $rng.End - 1
    $rng.Text = $NewText
}

# Row 4 = "Sprint-3"
Set-CellText $sprintTable 4 4 "06 Nov 2022"
Set-CellText $sprintTable 4 5 "11 Nov 2022"
Set-CellText $sprintTable 4 7 "09 Nov 2022"

# Row 5 = "Sprint-4"
Set-CellText $sprintTable 5 4 "10 Nov 2022"
Set-CellText $sprintTable 5 5 "15 Nov 2022"
Set-CellText $sprintTable 5 7 "15 Nov 2022"
